{"js": "// Update the date line above the table.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(\"2025-07-04 Friday\", Word.InsertLocation.replace);\n\n// Update all 100 math-problem cells in the 20x5 table in one shot.\n// newValues[row][col] lines up with the table in row-major reading order,\n// matching the order of the replacements in the diff.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\nconst newValues = [\n    [\"67-49=\", \"38+40=\", \"47+15=\", \"76+4=\", \"58+2=\"],\n    [\"4+46=\", \"31-21=\", \"50+21=\", \"78-67=\", \"55+28=\"],\n    [\"11+71=\", \"44-14=\", \"36+56=\", \"4+6=\", \"36+44=\"],\n    [\"89-62=\", \"72-57=\", \"24+32=\", \"93-19=\", \"56-38=\"],\n    [\"29+64=\", \"30+10=\", \"37+55=\", \"66+29=\", \"97-50=\"],\n    [\"92-54=\", \"23+46=\", \"6+92=\", \"47+12=\", \"91-19=\"],\n    [\"81-62=\", \"10-1=\", \"26+64=\", \"86-4=\", \"91+5=\"],\n    [\"54-42=\", \"20+15=\", \"6+21=\", \"3+95=\", \"57-55=\"],\n    [\"67-7=\", \"88+0=\", \"0+56=\", \"53-38=\", \"97-64=\"],\n    [\"31-28=\", \"53+26=\", \"70-9=\", \"94+4=\", \"29-5=\"],\n    [\"98-0=\", \"95-40=\", \"66-38=\", \"69-54=\", \"76-19=\"],\n    [\"61+16=\", \"52-50=\", \"48-43=\", \"45-21=\", \"49-9=\"],\n    [\"90-1=\", \"45+20=\", \"10+86=\", \"82+10=\", \"86+6=\"],\n    [\"37+1=\", \"95-88=\", \"71-57=\", \"46+37=\", \"82-4=\"],\n    [\"97-25=\", \"96-30=\", \"22+7=\", \"54+20=\", \"80-67=\"],\n    [\"9-4=\", \"44+33=\", \"91-31=\", \"8+68=\", \"17+43=\"],\n    [\"52-32=\", \"16+30=\", \"42+2=\", \"78-24=\", \"28+59=\"],\n    [\"48-32=\", \"50-45=\", \"75+18=\", \"80-13=\", \"63+22=\"],\n    [\"10+68=\", \"22+8=\", \"47-0=\", \"75-68=\", \"69-32=\"],\n    [\"42-29=\", \"85-1=\", \"53+23=\", \"65+27=\", \"28+23=\"]\n];\ntable.values = newValues;\nawait context.sync();\n\n", "ps1": "$d = $word.ActiveDocument\n\n# Paragraph 1 holds the date line above the table.\n$d.Paragraphs(1).Range.Text = \"2025-07-04 Friday\"\n\n# The 5x20 table of math problems: new values in row-major order\n# (row 1 col 1..5, row 2 col 1..5, ...), matching the order of the diff.\n$newValues = @(\n    \"67-49=\",\n    \"38+40=\",\n    \"47+15=\",\n    \"76+4=\",\n    \"58+2=\",\n    \"4+46=\",\n    \"31-21=\",\n    \"50+21=\",\n    \"78-67=\",\n    \"55+28=\",\n    \"11+71=\",\n    \"44-14=\",\n    \"36+56=\",\n    \"4+6=\",\n    \"36+44=\",\n    \"89-62=\",\n    \"72-57=\",\n    \"24+32=\",\n    \"93-19=\",\n    \"56-38=\",\n    \"29+64=\",\n    \"30+10=\",\n    \"37+55=\",\n    \"66+29=\",\n    \"97-50=\",\n    \"92-54=\",\n    \"23+46=\",\n    \"6+92=\",\n    \"47+12=\",\n    \"91-19=\",\n    \"81-62=\",\n    \"10-1=\",\n    \"26+64=\",\n    \"86-4=\",\n    \"91+5=\",\n    \"54-42=\",\n    \"20+15=\",\n    \"6+21=\",\n    \"3+95=\",\n    \"57-55=\",\n    \"67-7=\",\n    \"88+0=\",\n    \"0+56=\",\n    \"53-38=\",\n    \"97-64=\",\n    \"31-28=\",\n    \"53+26=\",\n    \"70-9=\",\n    \"94+4=\",\n    \"29-5=\",\n    \"98-0=\",\n    \"95-40=\",\n    \"66-38=\",\n    \"69-54=\",\n    \"76-19=\",\n    \"61+16=\",\n    \"52-50=\",\n    \"48-43=\",\n    \"45-21=\",\n    \"49-9=\",\n    \"90-1=\",\n    \"45+20=\",\n    \"10+86=\",\n    \"82+10=\",\n    \"86+6=\",\n    \"37+1=\",\n    \"95-88=\",\n    \"71-57=\",\n    \"46+37=\",\n    \"82-4=\",\n    \"97-25=\",\n    \"96-30=\",\n    \"22+7=\",\n    \"54+20=\",\n    \"80-67=\",\n    \"9-4=\",\n    \"44+33=\",\n    \"91-31=\",\n    \"8+68=\",\n    \"17+43=\",\n    \"52-32=\",\n    \"16+30=\",\n    \"42+2=\",\n    \"78-24=\",\n    \"28+59=\",\n    \"48-32=\",\n    \"50-45=\",\n    \"75+18=\",\n    \"80-13=\",\n    \"63+22=\",\n    \"10+68=\",\n    \"22+8=\",\n    \"47-0=\",\n    \"75-68=\",\n    \"69-32=\",\n    \"42-29=\",\n    \"85-1=\",\n    \"53+23=\",\n    \"65+27=\",\n    \"28+23=\"\n)\n\n$t = $d.Tables(1)\n$cols = $t.Columns.Count\n$idx = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $t.Cell($r, $c).Range.Text = $newValues[$idx]\n        $idx = $idx + 1\n    }\n}\n\n"}
